$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8: Average ride duration (Easy)
$ws.Range("C8").Value = "https://www.interviewquery.com/questions/average-ride-duration"
$ws.Range("A8").Value = "Average ride duration"
$ws.Range("B8").Value = "Easy"
$ws.Range("D8").Value = "for mysql, use timestampdiff or datediff for difference, for postgresql, use date_part or extract"

# Row 9: Top 5 turnover risk (Medium)
$ws.Range("C9").Value = "https://www.interviewquery.com/questions/top-5-turnover-risk"
$ws.Range("A9").Value = "Top 5 turnover risk"
$ws.Range("B9").Value = "Medium"
$ws.Range("D9").Value = "Solved using left join and subquery, but can be solved using inner join on first table and conditions"

# Row 10: Exam scores (Medium)
$ws.Range("C10").Value = "https://www.interviewquery.com/questions/exam-scores"
$ws.Range("A10").Value = "Exam scores"
$ws.Range("B10").Value = "Medium"
$ws.Range("D10").Value = "Tricky - https://learnsql.com/blog/case-when-with-sum/ - Can be solved with CASE WHEN THEN AND END or IF Condition, but since we are grouping by student names, an aggregate function has to exsit. In this case SUM; MAX can also work"

# Apply wrap-text style (style index 3) to the new C/D cells, matching existing rows
$ws.Range("C8:D10").WrapText = $true

# Row heights grow to fit the wrapped comment text, same as the other rows in the sheet
$ws.Rows.Item(8).RowHeight = 51
$ws.Rows.Item(9).RowHeight = 51
$ws.Rows.Item(10).RowHeight = 102

$ws.Range("D17").Select()
